# DemoQA_RegistrationSuite.xlsx update
# - A2 label changes from "RegistrationPage" to "RegistrationSuite"
# - Row 3 (previously "CustomerSuite" / "N") is cleared out
# - Active selection moves from A2 to C6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the suite referenced in row 2 from RegistrationPage -> RegistrationSuite
$ws.Range("A2").Value = "RegistrationSuite"

# Clear out the now-unused CustomerSuite / N row
$ws.Range("A3:B3").ClearContents() | Out-Null

# Move the active selection to C6
$ws.Range("C6").Select() | Out-Null
